# Extraction_MasterList.xlsx update
# - Buffer reagent renamed from "EB" to "C6" across existing extraction logs
# - New sheet "20_Nov_2025" added with a trial extraction for testing
#   purification kits on two test-soil samples (T1, T2)

$wb = $excel.ActiveWorkbook

# --- 1) Rename Buffer value EB -> C6 on the two existing sheets ---------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Replace("EB", "C6", 2)

$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Replace("EB", "C6", 2)

# --- 2) Add the new "20_Nov_2025" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "20_Nov_2025"

# Header row (same layout as the other extraction logs)
$ws3.Range("A1").Value = "Sample_ID"
$ws3.Range("B1").Value = "Label"
$ws3.Range("C1").Value = "Extraction_Date"
$ws3.Range("D1").Value = "Extraction_Num"
$ws3.Range("E1").Value = "Mass_g"
$ws3.Range("F1").Value = "Buffer"
$ws3.Range("G1").Value = "Conc_ngul"
$ws3.Range("H1").Value = "260_280"
$ws3.Range("I1").Value = "260_230"
$ws3.Range("J1").Value = "Vol_ul"
$ws3.Range("K1").Value = "Notes"

# Data rows - trial extraction on test soils (filled column-by-column so
# new shared-string entries land in the same order the author typed them:
# T1, T2, Test_Soils, For testing purification kits)
$ws3.Range("A2").Value = "T1"
$ws3.Range("A3").Value = "T2"

$ws3.Range("B2").Value = "Test_Soils "
$ws3.Range("B3").Value = "Test_Soils "

$ws3.Range("C2").Value = 45981
$ws3.Range("C3").Value = 45981

$ws3.Range("D2").Value = 1
$ws3.Range("D3").Value = 2

$ws3.Range("E2").Value = 0.35
$ws3.Range("E3").Value = 0.52

$ws3.Range("F2").Value = "C6"
$ws3.Range("F3").Value = "C6"

$ws3.Range("G2").Value = 53.2
$ws3.Range("G3").Value = 45.7

$ws3.Range("H2").Value = 1.87
$ws3.Range("H3").Value = 1.92

$ws3.Range("I2").Value = 0.85
$ws3.Range("I3").Value = 0.25

$ws3.Range("J2").Value = 100
$ws3.Range("J3").Value = 100

$ws3.Range("K2").Value = "For testing purification kits"
$ws3.Range("K3").Value = "For testing purification kits"

# --- 3) Restore / set the selections the author left on each tab --------
$ws1.Range("I15").Select()
$ws2.Range("A54:K56").Select()
$ws3.Select()
$ws3.Range("L4").Select()
